$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "code"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "descr"
$ws.Range("D1").Value = "lang_code"
$ws.Range("E1").Value = "is_active"
$ws.Range("F1").Value = "cr_by"
$ws.Range("G1").Value = "cr_dtimes"
$ws.Range("H1").Value = "upd_by"
$ws.Range("I1").Value = "upd_dtimes"
$ws.Range("J1").Value = "is_deleted"
$ws.Range("K1").Value = "del_dtimes"

# --- Common audit values ---
$crDate = 45079.577151851852

# --- Row 2 : DKS / Ordinateur ---
$ws.Range("A2").Value = "DKS"
$ws.Range("B2").Value = "Ordinateur"
$ws.Range("C2").Value = "Ordinateurs de bureau"
$ws.Range("D2").Value = "fra"
$ws.Range("E2").Value = $true
$ws.Range("F2").Value = "superadmin"
$ws.Range("G2").Value = $crDate
$ws.Range("H2").Value = "NULL"
$ws.Range("I2").Value = "NULL"
$ws.Range("J2").Value = $false
$ws.Range("K2").Value = "NULL"

# --- Row 3 : LTP / Portable ---
$ws.Range("A3").Value = "LTP"
$ws.Range("B3").Value = "Portable"
$ws.Range("C3").Value = "Ordinateurs portable"
$ws.Range("D3").Value = "fra"
$ws.Range("E3").Value = $true
$ws.Range("F3").Value = "superadmin"
$ws.Range("G3").Value = $crDate
$ws.Range("H3").Value = "NULL"
$ws.Range("I3").Value = "NULL"
$ws.Range("J3").Value = $false
$ws.Range("K3").Value = "NULL"

# --- Row 4 : TBT / Tablette ---
$ws.Range("A4").Value = "TBT"
$ws.Range("B4").Value = "Tablette"
$ws.Range("C4").Value = "Tablette"
$ws.Range("D4").Value = "fra"
$ws.Range("E4").Value = $true
$ws.Range("F4").Value = "superadmin"
$ws.Range("G4").Value = $crDate
$ws.Range("H4").Value = "NULL"
$ws.Range("I4").Value = "NULL"
$ws.Range("J4").Value = $false
$ws.Range("K4").Value = "NULL"

# --- Format the date/time columns like the source workbook (builtin mm:ss.0 format, id 47) ---
$ws.Range("G2:G4").NumberFormat = "mm:ss.0"

# --- Selection in the saved view matches the source file ---
$ws.Range("D11").Select() | Out-Null
